$wb = $excel.ActiveWorkbook

# Sheet ALC, row 17: One for the Road
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2980.024
$ws.Range("I17").Value = 295
$ws.Range("J17").Value = 3114.275
$ws.Range("K17").Value = 885
$ws.Range("L17").Value = 9342.825000000001
$ws.Range("M17").Value = -717
$ws.Range("N17").Value = -9678.825000000001

# Sheet ALC, row 19: Unbreak My Heart
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 33334442
$ws.Range("I19").Value = 111111510
$ws.Range("J19").Value = 1413.1428
$ws.Range("K19").Value = 111111510
$ws.Range("L19").Value = 1413.1428
$ws.Range("M19").Value = -111111335
$ws.Range("N19").Value = -1763.1428

# Sheet ALC, row 40: Stuck in the Moment
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1350.65
$ws.Range("I40").Value = 1445.3636
$ws.Range("J40").Value = 1234.8889
$ws.Range("K40").Value = 1445.3636
$ws.Range("L40").Value = 1234.8889
$ws.Range("M40").Value = -1270.3636
$ws.Range("N40").Value = -1584.8889

# Sheet ALC, row 113: Amaro Kart
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2702.7693
$ws.Range("I113").Value = 2125
$ws.Range("J113").Value = 2959.5557
$ws.Range("K113").Value = 2125
$ws.Range("L113").Value = 2959.5557
$ws.Range("M113").Value = 1129
$ws.Range("N113").Value = -9467.555700000001

# Sheet ALC, row 132: Fast-forwarding Flora
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1383.6492
$ws.Range("I132").Value = 1494.625
$ws.Range("K132").Value = 4483.875
$ws.Range("M132").Value = -1953.875

# Sheet ALC, row 137: Cutting Edge of Culinary Quality
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2322.116
$ws.Range("I137").Value = 1784.1632
$ws.Range("J137").Value = 3640.1
$ws.Range("K137").Value = 5352.4896
$ws.Range("L137").Value = 10920.3
$ws.Range("M137").Value = -2802.4896
$ws.Range("N137").Value = -16020.3

# Sheet ALC, row 138: All-night Crafting
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2634.2856
$ws.Range("I138").Value = 1336.1818
$ws.Range("J138").Value = 4496.7827
$ws.Range("K138").Value = 4008.5454
$ws.Range("L138").Value = 13490.3481
$ws.Range("M138").Value = 1131.4546
$ws.Range("N138").Value = -23770.3481

# Sheet ALC, row 141: Remedy for Reason
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3170.6858
$ws.Range("I141").Value = 2161.2083
$ws.Range("J141").Value = 5373.1816
$ws.Range("K141").Value = 6483.624899999999
$ws.Range("L141").Value = 16119.5448
$ws.Range("M141").Value = -1303.624899999999
$ws.Range("N141").Value = -26479.5448

# Sheet ARM, row 45: Hollow Hallmarks
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1424.8334
$ws.Range("I45").Value = 1397.36
$ws.Range("K45").Value = 1397.36
$ws.Range("M45").Value = -1020.36

# Sheet ARM, row 110: Scheduled Maintenance
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1442.3572
$ws.Range("I110").Value = 1455.5454
$ws.Range("J110").Value = 1394
$ws.Range("K110").Value = 1455.5454
$ws.Range("L110").Value = 1394
$ws.Range("M110").Value = 589.4546
$ws.Range("N110").Value = -5484

# Sheet BSM, row 99: Meddle in Metal
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1088.8889
$ws.Range("I99").Value = 1230
$ws.Range("J99").Value = 999.0909
$ws.Range("K99").Value = 1230
$ws.Range("L99").Value = 999.0909
$ws.Range("M99").Value = 268
$ws.Range("N99").Value = -3995.0909

# Sheet BSM, row 131: Plying with Precision
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H131").Value = 35000
$ws.Range("J131").Value = 35000
$ws.Range("L131").Value = 35000
$ws.Range("N131").Value = -45080

# Sheet BSM, row 134: Ruthenium Supremium
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 36346.9
$ws.Range("I134").Value = 3330.4736
$ws.Range("J134").Value = 93375.27
$ws.Range("K134").Value = 9991.4208
$ws.Range("L134").Value = 280125.81
$ws.Range("M134").Value = -7456.4208
$ws.Range("N134").Value = -285195.81

# Sheet CRP, row 10: Spears and Sorcery
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 36002.668
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 36002.668
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 36002.668
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -36280.668

# Sheet CRP, row 31: Wall Not Found
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6201.533
$ws.Range("I31").Value = 7116
$ws.Range("K31").Value = 7116
$ws.Range("M31").Value = -6821

# Sheet CRP, row 34: Armoires of the Rich and Famous
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6201.533
$ws.Range("I34").Value = 7116
$ws.Range("K34").Value = 7116
$ws.Range("M34").Value = -6914

# Sheet CRP, row 99: O Pine
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3370.4707
$ws.Range("I99").Value = 2454
$ws.Range("J99").Value = 5570
$ws.Range("K99").Value = 2454
$ws.Range("L99").Value = 5570
$ws.Range("M99").Value = -956
$ws.Range("N99").Value = -8566

# Sheet CRP, row 126: A Better Conductor
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3370.4707
$ws.Range("I126").Value = 2454
$ws.Range("J126").Value = 5570
$ws.Range("K126").Value = 7362
$ws.Range("L126").Value = 16710
$ws.Range("M126").Value = -4892
$ws.Range("N126").Value = -21650

# Sheet CUL, row 114: One Last Meal
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 693.6667
$ws.Range("J114").Value = 1031
$ws.Range("L114").Value = 3093
$ws.Range("N114").Value = -9601

# Sheet CUL, row 129: Comfort Food
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2436.5264
$ws.Range("J129").Value = 1950.4546
$ws.Range("L129").Value = 5851.3638
$ws.Range("N129").Value = -15851.3638

# Sheet CUL, row 131: The Mountain Steeped
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 28497.97
$ws.Range("I131").Value = 2054.4443
$ws.Range("J131").Value = 37651.5
$ws.Range("K131").Value = 6163.3329
$ws.Range("L131").Value = 112954.5
$ws.Range("M131").Value = -1123.3329
$ws.Range("N131").Value = -123034.5

# Sheet CUL, row 134: Don't Knock It Till You've Tried It
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 3738.6775
$ws.Range("I134").Value = 3362.9565
$ws.Range("J134").Value = 4818.875
$ws.Range("K134").Value = 10088.8695
$ws.Range("L134").Value = 14456.625
$ws.Range("M134").Value = -5018.869499999999
$ws.Range("N134").Value = -24596.625

# Sheet CUL, row 137: Creative Chocolate
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 15495.415
$ws.Range("I137").Value = 1315.45
$ws.Range("J137").Value = 29000.143
$ws.Range("K137").Value = 3946.35
$ws.Range("L137").Value = 87000.429
$ws.Range("M137").Value = 1153.65
$ws.Range("N137").Value = -97200.429

# Sheet GSM, row 11: A Ringing Success
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 366666660
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

# Sheet GSM, row 57: Gold Is So Last Year
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 18000
$ws.Range("J57").Value = 18000
$ws.Range("L57").Value = 18000
$ws.Range("N57").Value = -19640

# Sheet GSM, row 62: The Goggles, They Do Naught
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 49000
$ws.Range("J62").Value = 49000
$ws.Range("L62").Value = 49000
$ws.Range("N62").Value = -50372

# Sheet GSM, row 65: Peril Never Wore Safety Goggles (L)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H65").Value = 49000
$ws.Range("J65").Value = 49000
$ws.Range("L65").Value = 147000
$ws.Range("N65").Value = -153864

# Sheet GSM, row 102: Put the Metal to the Peddle
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3631.3333
$ws.Range("I102").Value = 3728.762
$ws.Range("J102").Value = 3460.8333
$ws.Range("K102").Value = 3728.762
$ws.Range("L102").Value = 3460.8333
$ws.Range("M102").Value = -2106.762
$ws.Range("N102").Value = -6704.8333

# Sheet LTW, row 13: Throwing Down the Gauntlet
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

# Sheet LTW, row 55: It's Not a Job, It's a Calling
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 330.1579
$ws.Range("I55").Value = 247.9
$ws.Range("J55").Value = 421.55554
$ws.Range("K55").Value = 247.9
$ws.Range("L55").Value = 421.55554
$ws.Range("M55").Value = -74.90000000000001
$ws.Range("N55").Value = -767.5555400000001

# Sheet LTW, row 63: From Mud to Mourning
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 41723.332
$ws.Range("J63").Value = 41723.332
$ws.Range("L63").Value = 41723.332
$ws.Range("N63").Value = -43221.332

# Sheet LTW, row 66: These Boots Are Made for Hawkin' (L)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H66").Value = 41723.332
$ws.Range("J66").Value = 41723.332
$ws.Range("L66").Value = 125169.996
$ws.Range("N66").Value = -132657.996

# Sheet WVR, row 43: Walk Softly and Carry a Big Halberd
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 24251.4
$ws.Range("I43").Value = 10027
$ws.Range("J43").Value = 27807.5
$ws.Range("K43").Value = 10027
$ws.Range("L43").Value = 27807.5
$ws.Range("M43").Value = -9878
$ws.Range("N43").Value = -28105.5

# Sheet WVR, row 58: Seeing It Through to the End
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 18000
$ws.Range("J58").Value = 18000
$ws.Range("L58").Value = 18000
$ws.Range("N58").Value = -18616

# Sheet WVR, row 122: Heavy Armoire
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1809.2727
$ws.Range("I122").Value = 1362.75
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 4088.25
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -1638.25
$ws.Range("N122").Value = -13900

# Sheet WVR, row 137: Traditional Trousers
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 45000
$ws.Range("J137").Value = 60000
$ws.Range("L137").Value = 60000
$ws.Range("N137").Value = -70200
